$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(18)
$tr = $shape.TextFrame.TextRange

# Merge the three runs that make up "Retour ... majeure" into a single run
$run1 = $tr.Characters(82, 56)
$run1.Text = "Retour à « pour implementation » si modification majeure"

# Merge the "P" + "assage" runs into a single "Passage" run
$run2 = $tr.Characters(222, 7)
$run2.Text = "Passage"
